$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.898.85"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.762.42"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.03"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.14"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "3.761.56"
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("E11").Value = "  -5.20%  "
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.38"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "4.377.37"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "3.762.33"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "69.946.00"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.78"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "509.22"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.39"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.17"
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.16"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  +20.55%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.18"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("E36").Value = "  +4.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.21"
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.13"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.72"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.89"
$ws.Range("E43").Value = "  +2.93%  "
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").Value = "3.008.14"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").Value = "  -3.51%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.65"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  +2.27%  "
